$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Abfragen" row score cell: the text is split across two runs
#    ("……. / 5" + " Punkten") with a leftover "_GoBack" bookmark sitting at
#    the split point. Replacing the phrase in a range around that bookmark
#    collapses the two runs back into one and drops the bookmark markup.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$winStart = [Math]::Max(0, $bm.Start - 20)
$fixRange = $d.Range($winStart, $bm.Start + 20)
$fixRange.Find.Execute("……. / 5 Punkten", $false, $false, $false, $false, `
    $false, $true, 1, $false, "……. / 5 Punkten", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Abfrage: {hw_9_file}" heading: split the run after "…hw_9_fil" and plant
#    a fresh "_GoBack" bookmark at the cut, mirroring where it now lives.
# ---------------------------------------------------------------------------
$hdr = $d.Content
$hdr.Find.Execute("Abfrage: {hw_9_file}", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $hdr.Start + "Abfrage: {hw_9_fil".Length
$cut = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $cut) | Out-Null

# ---------------------------------------------------------------------------
# 3) Drop the four trailing empty paragraphs after the "{@hw_9}" merge field,
#    right before the final sectPr.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "{@hw_9}") {
        $lastKeep = $p
    }
}
$cleanup = $d.Range($lastKeep.Range.End, $d.Content.End)
$cleanup.Delete()
